$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175829768180847
$ws.Range("B1").Value = 2.408564567565918
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.343636751174927
$ws.Range("E1").Value = 1.206317543983459
